$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, pushing the existing rows 11-14 down to 12-15.
# Excel carries formulas/number-formats/ranges along automatically.
$ws.Rows.Item(11).Insert()

# --- Populate the new battery row (row 11) ---
$ws.Range("A11").Value = "BATT"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "CR2032"
$ws.Range("E11").Value = "CR2032"
$ws.Range("F11").Value = 0.28999999999999998
$ws.Range("G11").Formula = "=F11*B11"
$ws.Range("H11").Value = 0.15734000000000001
$ws.Range("I11").Formula = "=H11*B11"
$ws.Range("J11").Value = "https://www.digikey.com/product-detail/en/panasonic-bsg/CR2032/P189-ND/31939"

# Copy the look (styles) of the row above onto the new row's lettered cells
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("F10:I10").Copy() | Out-Null
$ws.Range("F11:I11").PasteSpecial(-4122) | Out-Null

# The old row 11 (now row 12) loses its D-column "PIN" label in the new layout
$ws.Range("D12").Clear() | Out-Null

# --- Row 6: record the connector's MPN in column E ---
$ws.Range("E6").Value = "PRPC003DAAN-RC"

# --- Rows 7/8: the Part Page links were swapped ---
$ws.Range("J7").Value = "https://www.digikey.com/product-detail/en/yageo/RC0805FR-0710KL/311-10.0KCRCT-ND/730482"
$ws.Range("J8").Value = "https://www.digikey.com/product-detail/en/yageo/RC0805FR-0782RL/311-82.0CRCT-ND/731111"

# --- Selection / view state ---
$ws.Range("D11").Select() | Out-Null

$ws.Calculate()
